$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D, E
$ws.Range('D2').Value = '45.231.65'
$ws.Range('E2').Value = '  +4.81%  '

# Row 3: D, E
$ws.Range('D3').Value = '2.355.62'
$ws.Range('E3').Value = '  +1.62%  '

# Row 4: E
$ws.Range('E4').Value = '  -0.05%  '

# Row 5: D, E
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '109.20'
$ws.Range("D5").Style = "Normal"
$ws.Range('E5').Value = '  +0.79%  '

# Row 6: D, E
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '308.17'
$ws.Range("D6").Style = "Normal"
$ws.Range('E6').Value = '  -1.16%  '

# Row 7: D, E
$ws.Range("D7").NumberFormat = "@"
$ws.Range('D7').Value = '0.627'
$ws.Range("D7").Style = "Normal"
$ws.Range('E7').Value = '  -0.05%  '

# Row 8: E
$ws.Range('E8').Value = '  -0.15%  '

# Row 9: D, E
$ws.Range("D9").NumberFormat = "@"
$ws.Range('D9').Value = '0.614'
$ws.Range("D9").Style = "Normal"
$ws.Range('E9').Value = '  +1.01%  '

# Row 10: D, E
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '41.06'
$ws.Range("D10").Style = "Normal"
$ws.Range('E10').Value = '  +1.78%  '

# Row 11: D, E
$ws.Range("D11").NumberFormat = "@"
$ws.Range('D11').Value = '0.0914'
$ws.Range("D11").Style = "Normal"
$ws.Range('E11').Value = '  -0.10%  '

# Row 12: D, E
$ws.Range("D12").NumberFormat = "@"
$ws.Range('D12').Value = '8.41'
$ws.Range("D12").Style = "Normal"
$ws.Range('E12').Value = '  -0.04%  '

# Row 13: E
$ws.Range('E13').Value = '  +1.17%  '

# Row 14: D, E
$ws.Range("D14").NumberFormat = "@"
$ws.Range('D14').Value = '0.982'
$ws.Range("D14").Style = "Normal"
$ws.Range('E14').Value = '  -1.77%  '

# Row 15: D, E
$ws.Range('D15').Value = '2.714.52'
$ws.Range('E15').Value = '  +1.71%  '

# Row 16: D, E
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '15.32'
$ws.Range("D16").Style = "Normal"
$ws.Range('E16').Value = '  -0.71%  '

# Row 17: D, E
$ws.Range('D17').Value = '2.350.17'
$ws.Range('E17').Value = '  +1.50%  '

# Row 18: D, E
$ws.Range('D18').Value = '45.161.85'
$ws.Range('E18').Value = '  +5.42%  '

# Row 19: D, E
$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '7.26'
$ws.Range("D19").Style = "Normal"
$ws.Range('E19').Value = '  -3.54%  '

# Row 20: D, E
$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '0.0000106'
$ws.Range("D20").Style = "Normal"
$ws.Range('E20').Value = '  +0.56%  '

# Row 21: D, E
$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '13.08'
$ws.Range("D21").Style = "Normal"
$ws.Range('E21').Value = '  -0.74%  '

# Row 22: D, E
$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '73.24'
$ws.Range("D22").Style = "Normal"
$ws.Range('E22').Value = '  -0.70%  '

# Row 23: D, E
$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '3.42'
$ws.Range("D23").Style = "Normal"
$ws.Range('E23').Value = '  -2.03%  '

# Row 24: D, E
$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '258.61'
$ws.Range("D24").Style = "Normal"
$ws.Range('E24').Value = '  -3.10%  '

# Row 25: E
$ws.Range('E25').Value = '  +0.70%  '

# Row 26: E
$ws.Range('E26').Value = '  -0.26%  '

# Row 27: D, E
$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '11.05'
$ws.Range("D27").Style = "Normal"
$ws.Range('E27').Value = '  +0.28%  '

# Row 28: D, E
$ws.Range("D28").NumberFormat = "@"
$ws.Range('D28').Value = '7.33'
$ws.Range("D28").Style = "Normal"
$ws.Range('E28').Value = '  -5.79%  '

# Row 29: E
$ws.Range('E29').Value = '  +2.29%  '

# Row 30: D, E
$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '0.0962'
$ws.Range("D30").Style = "Normal"
$ws.Range('E30').Value = '  +10.32%  '

# Row 31: D, E
$ws.Range("D31").NumberFormat = "@"
$ws.Range('D31').Value = '22.32'
$ws.Range("D31").Style = "Normal"
$ws.Range('E31').Value = '  -0.47%  '

# Row 32: D, E
$ws.Range("D32").NumberFormat = "@"
$ws.Range('D32').Value = '37.70'
$ws.Range("D32").Style = "Normal"
$ws.Range('E32').Value = '  -2.77%  '

# Row 33: D, E
$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '170.69'
$ws.Range("D33").Style = "Normal"
$ws.Range('E33').Value = '  +2.78%  '

# Row 34: E
$ws.Range('E34').Value = '  +6.71%  '

# Row 35: E
$ws.Range('E35').Value = '  +0.05%  '

# Row 36: D, E
$ws.Range("D36").NumberFormat = "@"
$ws.Range('D36').Value = '4.81'
$ws.Range("D36").Style = "Normal"
$ws.Range('E36').Value = '  +3.24%  '

# Row 37: D, E
$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '0.114'
$ws.Range("D37").Style = "Normal"
$ws.Range('E37').Value = '  +1.30%  '

# Row 38: D, E
$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '2.95'
$ws.Range("D38").Style = "Normal"
$ws.Range('E38').Value = '  +4.88%  '

# Row 39: D, E
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '0.0356'
$ws.Range("D39").Style = "Normal"
$ws.Range('E39').Value = '  -0.25%  '

# Row 40: D, E
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '3.89'
$ws.Range("D40").Style = "Normal"
$ws.Range('E40').Value = '  +6.25%  '

# Row 41: E
$ws.Range('E41').Value = '  +7.30%  '

# Row 42: D, E
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '99.02'
$ws.Range("D42").Style = "Normal"
$ws.Range('E42').Value = '  -5.70%  '

# Row 43: D, E
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '0.231'
$ws.Range("D43").Style = "Normal"
$ws.Range('E43').Value = '  -0.40%  '

# Row 44: D, E
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '69.51'
$ws.Range("D44").Style = "Normal"
$ws.Range('E44').Value = '  -2.54%  '

# Row 45: B, C, D, E
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range('E45').Value = '  +0.21%  '

# Row 46: B, C, D, E
$ws.Range('B46').Value = 'Celestia'
$ws.Range('C46').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '12.71'
$ws.Range("D46").Style = "Normal"
$ws.Range('E46').Value = '  +2.07%  '

# Row 47: D, E
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '81.37'
$ws.Range("D47").Style = "Normal"
$ws.Range('E47').Value = '  +5.41%  '

# Row 48: B, C, D, E
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '112.36'
$ws.Range("D48").Style = "Normal"
$ws.Range('E48').Value = '  -0.74%  '

# Row 49: B, C, D, E
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '9.30'
$ws.Range("D49").Style = "Normal"
$ws.Range('E49').Value = '  +4.70%  '

# Row 50: D, E
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '5.49'
$ws.Range("D50").Style = "Normal"
$ws.Range('E50').Value = '  +4.47%  '

# Row 51: D, E
$ws.Range('D51').Value = '1.626.85'
$ws.Range('E51').Value = '  -3.40%  '
